$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrease the "farms_total_count" (C=4) and "farms_to_examine_count" (C=5)
# values by 8 for every existing weekly block (rows 2..121, in groups of 5).
for ($row = 2; $row -le 117; $row += 5) {
    $cell1 = $ws.Cells.Item($row, 4)
    $cell1.Value = $cell1.Value() - 8

    $cell2 = $ws.Cells.Item($row + 1, 4)
    $cell2.Value = $cell2.Value() - 8
}

# Append a new weekly block of data (YearWeekIso 202512 / LastDayOfWeek 2025-03-23).
$newRows = @(
    @{ C = "farms_total_count";            D = 12514 },
    @{ C = "farms_to_examine_count";        D = 1746 },
    @{ C = "farms_examined_count";          D = 10768 },
    @{ C = "farms_examined_positive_count"; D = 1468 },
    @{ C = "farms_examined_negative_count"; D = 9300 }
)

$r = 122
foreach ($item in $newRows) {
    $ws.Cells.Item($r, 1).Value = 202512

    # Copy the date cell's formatting (style) from the previous block, then
    # overwrite the value so it keeps the existing date number format.
    $srcDate = $ws.Cells.Item($r - 5, 2)
    $dstDate = $ws.Cells.Item($r, 2)
    $srcDate.Copy($dstDate)
    $dstDate.Value = 45739

    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $r++
}

# Update the active cell selection to match the author's recorded state.
$ws.Range("C8").Select()
